$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value for every data row (2..532).
# All of these cells currently store the Excel serial date 45178 (2023-09-09)
# and need to be bumped by one day to 45179 (2023-09-10).
$range = $ws.Range("C2:C532")
$range.Value = 45179
